$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the "days elapsed" formula down column I (I2:I378), referencing the
# matching row of column G (=TODAY()-G<row>). Writing it in the same
# chunk sizes Excel itself would use reproduces the shared-formula grouping
# (si=0..5) that appears in the saved workbook.
$ws.Range("I2").Formula        = "=TODAY()-G2"
$ws.Range("I3:I66").Formula    = "=TODAY()-G3"
$ws.Range("I67:I130").Formula  = "=TODAY()-G67"
$ws.Range("I131:I194").Formula = "=TODAY()-G131"
$ws.Range("I195:I258").Formula = "=TODAY()-G195"
$ws.Range("I259:I322").Formula = "=TODAY()-G259"
$ws.Range("I323:I378").Formula = "=TODAY()-G323"

# Reproduce the resulting selection/scroll position left behind in the file:
# the whole filled range is selected with I2 as the active cell, and the
# window has been scrolled down so row 190 is at the top.
$ws.Range("I2:I378").Select()
$excel.ActiveWindow.ScrollRow = 190
